# "getting the last version" - refresh the sheet's last-used-row counter
# and move the saved cursor/scroll position to the bottom-right corner of
# the used range (column Y, row 2), matching the last-viewed cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the tracked counter in Y2 (434 -> 435).
$ws.Range("Y2").Value = 435

# Scroll the window so column V is the first visible column, then park the
# selection on the last cell (Y2) - this is what was active when the file
# was last saved.
$excel.ActiveWindow.ScrollColumn = 22
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Y2").Select()
